# edit.ps1 -- apply the tracked changes described by the commit diff:
#   1. Update the document date.
#   2. Fix the table layout (switch to fixed layout) on the two
#      percentage-width tables ("Pooled estimated values ..." and
#      "Linear regime parameter estimates ...").
#   3. Add a new "Abstract Title" paragraph style and tighten the
#      spacing above the "Abstract" style.
#   4. Add a new "Footnote Block Text" paragraph style (based on
#      "Footnote Text", mirrors the existing "Block Text" style).

$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1. Date: 2024-03-27 -> 2024-07-03
# -------------------------------------------------------------------
$d.Content.Find.Execute("2024-03-27", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-03", 2) | Out-Null

# -------------------------------------------------------------------
# 2. Force fixed table layout on the two pct-width tables that got
#    <w:tblLayout w:type="fixed"/> added to their tblPr.
# -------------------------------------------------------------------
foreach ($t in $d.Tables) {
    if ($t.PreferredWidthType -eq 2) {
        # wdPreferredWidthPercent == 2 -> these are the tblW type="pct" tables
        $t.AllowAutoFit = $false
    }
}

# -------------------------------------------------------------------
# 3a. New style: "Abstract Title"
# -------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# -------------------------------------------------------------------
# 3b. "Abstract" style: space-before 300 -> 100 (twentieths of a point
#     i.e. 15pt -> 5pt)
# -------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# -------------------------------------------------------------------
# 4. New style: "Footnote Block Text" (based on "Footnote Text")
# -------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "edit applied"
